$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CBM thickness")

# Row 7 (Lash et al., 1989 (11 wk.))
$ws.Range("C7").Comment.Delete()
$ws.Range("G7").Comment.Delete()

$ws.Range("B7").Value = 61.87
$ws.Range("C7").ClearContents()
$ws.Range("E7").Value = 1.33
$ws.Range("F7").Value = 68.13
$ws.Range("G7").ClearContents()
$ws.Range("I7").Value = 1.66

# Row 8 (Lash et al., 1989 (18 wk.))
$ws.Range("C8").Comment.Delete()
$ws.Range("G8").Comment.Delete()

$ws.Range("B8").Value = 55.67
$ws.Range("C8").ClearContents()
$ws.Range("E8").Value = 1.04
$ws.Range("F8").Value = 57.82
$ws.Range("G8").ClearContents()
$ws.Range("I8").Value = 1.24

$ws.Activate()
$ws.Range("I9").Select()
